$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 289.14285
$ws.Range("J9").Value = 335
$ws.Range("L9").Value = 335
$ws.Range("N9").Value = -673
$ws.Range("H141").Value = 1200
$ws.Range("I141").Value = 1200
$ws.Range("K141").Value = 3600
$ws.Range("M141").Value = 1580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3682.8333
$ws.Range("I74").Value = 4059.8
$ws.Range("K74").Value = 4059.8
$ws.Range("M74").Value = -3185.8
$ws.Range("H77").Value = 3682.8333
$ws.Range("I77").Value = 4059.8
$ws.Range("K77").Value = 20299
$ws.Range("M77").Value = -15931
$ws.Range("H132").Value = 2534.6
$ws.Range("I132").Value = 2144.2144
$ws.Range("K132").Value = 6432.6432
$ws.Range("M132").Value = -3902.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 193.47058
$ws.Range("I7").Value = 198.63637
$ws.Range("J7").Value = 184
$ws.Range("K7").Value = 198.63637
$ws.Range("L7").Value = 184
$ws.Range("M7").Value = -85.63637
$ws.Range("N7").Value = -410
$ws.Range("H31").Value = 8060.5454
$ws.Range("I31").Value = 7914.4116
$ws.Range("K31").Value = 7914.4116
$ws.Range("M31").Value = -7619.4116
$ws.Range("H34").Value = 8060.5454
$ws.Range("I34").Value = 7914.4116
$ws.Range("K34").Value = 7914.4116
$ws.Range("M34").Value = -7712.4116
$ws.Range("H58").Value = 2424
$ws.Range("I58").Value = 2424
$ws.Range("K58").Value = 2424
$ws.Range("M58").Value = -2221
$ws.Range("H99").Value = 759.6
$ws.Range("I99").Value = 699.5
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 699.5
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 798.5
$ws.Range("N99").Value = -3996
$ws.Range("H126").Value = 759.6
$ws.Range("I126").Value = 699.5
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2098.5
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 371.5
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 3203
$ws.Range("I132").Value = 3602.2856
$ws.Range("J132").Value = 2504.25
$ws.Range("K132").Value = 10806.8568
$ws.Range("L132").Value = 7512.75
$ws.Range("M132").Value = -8276.856800000001
$ws.Range("N132").Value = -12572.75
$ws.Range("H134").Value = 2249.5
$ws.Range("I134").Value = 2249.5
$ws.Range("K134").Value = 6748.5
$ws.Range("M134").Value = -4213.5
$ws.Range("H136").Value = 2424
$ws.Range("I136").Value = 2424
$ws.Range("K136").Value = 7272
$ws.Range("M136").Value = -4722

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 17
$ws.Range("I2").Value = 19.5
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 72
$ws.Range("M2").Value = -4
$ws.Range("N2").Value = -298
$ws.Range("H7").Value = 26
$ws.Range("I7").Value = 26
$ws.Range("J7").Value = 26
$ws.Range("K7").Value = 78
$ws.Range("L7").Value = 78
$ws.Range("M7").Value = 34
$ws.Range("N7").Value = -302
$ws.Range("H23").Value = 1307.5385
$ws.Range("I23").Value = 1470.5714
$ws.Range("J23").Value = 1117.3334
$ws.Range("K23").Value = 4411.7142
$ws.Range("L23").Value = 3352.0002
$ws.Range("M23").Value = -4176.7142
$ws.Range("N23").Value = -3822.0002
$ws.Range("H132").Value = 3099.8
$ws.Range("I132").Value = 3125
$ws.Range("K132").Value = 28125
$ws.Range("M132").Value = -25595

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 36.615383
$ws.Range("J2").Value = 50.333332
$ws.Range("L2").Value = 50.333332
$ws.Range("N2").Value = -276.333332
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H113").Value = 799.5
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 799
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 799
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -5139
$ws.Range("H126").Value = 4870.5
$ws.Range("I126").Value = 5160.6665
$ws.Range("K126").Value = 15481.9995
$ws.Range("M126").Value = -13011.9995
$ws.Range("H132").Value = 6338.5
$ws.Range("I132").Value = 6338.5
$ws.Range("K132").Value = 19015.5
$ws.Range("M132").Value = -16485.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 3828.3333
$ws.Range("I7").Value = 3034.1667
$ws.Range("J7").Value = 5416.6665
$ws.Range("K7").Value = 3034.1667
$ws.Range("L7").Value = 5416.6665
$ws.Range("M7").Value = -2922.1667
$ws.Range("N7").Value = -5640.6665
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H126").Value = 3828.3333
$ws.Range("I126").Value = 3034.1667
$ws.Range("J126").Value = 5416.6665
$ws.Range("K126").Value = 9102.500100000001
$ws.Range("L126").Value = 16249.9995
$ws.Range("M126").Value = -6632.500100000001
$ws.Range("N126").Value = -21189.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H100").Value = 5131
$ws.Range("I100").Value = 5339.533
$ws.Range("J100").Value = 2003
$ws.Range("K100").Value = 10679.066
$ws.Range("L100").Value = 4006
$ws.Range("M100").Value = -10138.066
$ws.Range("N100").Value = -5088
$ws.Range("H126").Value = 1349.5
$ws.Range("I126").Value = 1399
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 4197
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1727
$ws.Range("N126").Value = -8840
